$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H20").Value = 2131.5
$ws.Range("I20").Value = 2131.5
$ws.Range("K20").Value = 2131.5
$ws.Range("M20").Value = -1901.5

$ws.Range("H35").Value = 2131.5
$ws.Range("I35").Value = 2131.5
$ws.Range("K35").Value = 2131.5
$ws.Range("M35").Value = -1752.5

$ws.Range("H40").Value = 1685.7142
$ws.Range("I40").Value = 2800
$ws.Range("J40").Value = 1500
$ws.Range("K40").Value = 2800
$ws.Range("L40").Value = 1500
$ws.Range("M40").Value = -2625
$ws.Range("N40").Value = -1850

$ws.Range("H76").Value = 8674.947
$ws.Range("I76").Value = 12181.818
$ws.Range("J76").Value = 3853
$ws.Range("K76").Value = 12181.818
$ws.Range("L76").Value = 3853
$ws.Range("M76").Value = -11866.818
$ws.Range("N76").Value = -4483

$ws.Range("H79").Value = 8674.947
$ws.Range("I79").Value = 12181.818
$ws.Range("J79").Value = 3853
$ws.Range("K79").Value = 12181.818
$ws.Range("L79").Value = 3853
$ws.Range("M79").Value = -11089.818
$ws.Range("N79").Value = -6037

$ws.Range("H112").Value = 22728822
$ws.Range("I112").Value = 250000270
$ws.Range("J112").Value = 1676.7
$ws.Range("K112").Value = 750000810
$ws.Range("L112").Value = 5030.1
$ws.Range("M112").Value = -749999702
$ws.Range("N112").Value = -7246.1

$ws.Range("H129").Value = 837.1212
$ws.Range("I129").Value = 702.05554
$ws.Range("J129").Value = 999.2
$ws.Range("K129").Value = 2106.16662
$ws.Range("L129").Value = 2997.6
$ws.Range("M129").Value = 2893.83338
$ws.Range("N129").Value = -12997.6

$ws.Range("H132").Value = 1066497.6
$ws.Range("I132").Value = 1166.85
$ws.Range("J132").Value = 8168702.5
$ws.Range("K132").Value = 3500.55
$ws.Range("L132").Value = 24506107.5
$ws.Range("M132").Value = -970.5499999999997
$ws.Range("N132").Value = -24511167.5

$ws.Range("H134").Value = 30000
$ws.Range("J134").Value = 30000
$ws.Range("L134").Value = 30000
$ws.Range("N134").Value = -40140

$ws.Range("H138").Value = 2151.8545
$ws.Range("J138").Value = 1899.2222
$ws.Range("L138").Value = 5697.6666
$ws.Range("N138").Value = -15977.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1814.2916
$ws.Range("I2").Value = 1419.4
$ws.Range("J2").Value = 2472.4443
$ws.Range("K2").Value = 1419.4
$ws.Range("L2").Value = 2472.4443
$ws.Range("M2").Value = -1306.4
$ws.Range("N2").Value = -2698.4443

$ws.Range("H74").Value = 6148394
$ws.Range("I74").Value = 9655050
$ws.Range("J74").Value = 70190
$ws.Range("K74").Value = 9655050
$ws.Range("L74").Value = 70190
$ws.Range("M74").Value = -9654176
$ws.Range("N74").Value = -71938

$ws.Range("H77").Value = 6148394
$ws.Range("I77").Value = 9655050
$ws.Range("J77").Value = 70190
$ws.Range("K77").Value = 48275250
$ws.Range("L77").Value = 350950
$ws.Range("M77").Value = -48270882
$ws.Range("N77").Value = -359686

$ws.Range("H88").Value = 8492.817999999999
$ws.Range("I88").Value = 3995.3333
$ws.Range("J88").Value = 10179.375
$ws.Range("K88").Value = 3995.3333
$ws.Range("L88").Value = 10179.375
$ws.Range("M88").Value = -3589.3333
$ws.Range("N88").Value = -10991.375

$ws.Range("H91").Value = 8492.817999999999
$ws.Range("I91").Value = 3995.3333
$ws.Range("J91").Value = 10179.375
$ws.Range("K91").Value = 3995.3333
$ws.Range("L91").Value = 10179.375
$ws.Range("M91").Value = -2591.3333
$ws.Range("N91").Value = -12987.375

$ws.Range("H116").Value = 1814.2916
$ws.Range("I116").Value = 1419.4
$ws.Range("J116").Value = 2472.4443
$ws.Range("K116").Value = 1419.4
$ws.Range("L116").Value = 2472.4443
$ws.Range("M116").Value = 874.5999999999999
$ws.Range("N116").Value = -7060.4443

$ws.Range("H122").Value = 9261363
$ws.Range("I122").Value = 2324.1
$ws.Range("K122").Value = 6972.299999999999
$ws.Range("M122").Value = -4522.299999999999

$ws.Range("H132").Value = 108426.266
$ws.Range("I132").Value = 127287.625
$ws.Range("J132").Value = 94708.91
$ws.Range("K132").Value = 381862.875
$ws.Range("L132").Value = 284126.73
$ws.Range("M132").Value = -379332.875
$ws.Range("N132").Value = -289186.73

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1814.2916
$ws.Range("I3").Value = 1419.4
$ws.Range("J3").Value = 2472.4443
$ws.Range("K3").Value = 1419.4
$ws.Range("L3").Value = 2472.4443
$ws.Range("M3").Value = -1305.4
$ws.Range("N3").Value = -2700.4443

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 2678.2415
$ws.Range("I99").Value = 2318.524
$ws.Range("K99").Value = 2318.524
$ws.Range("M99").Value = -820.5239999999999

$ws.Range("H122").Value = 1477.9688
$ws.Range("I122").Value = 1192.4073
$ws.Range("J122").Value = 3020
$ws.Range("K122").Value = 3577.2219
$ws.Range("L122").Value = 9060
$ws.Range("M122").Value = -1127.2219
$ws.Range("N122").Value = -13960

$ws.Range("H126").Value = 2678.2415
$ws.Range("I126").Value = 2318.524
$ws.Range("K126").Value = 6955.572
$ws.Range("M126").Value = -4485.572

$ws.Range("H134").Value = 29891.775
$ws.Range("I134").Value = 1747.7188
$ws.Range("J134").Value = 142468
$ws.Range("K134").Value = 5243.1564
$ws.Range("L134").Value = 427404
$ws.Range("M134").Value = -2708.1564
$ws.Range("N134").Value = -432474

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H103").Value = 2461.7666
$ws.Range("I103").Value = 458.33334
$ws.Range("J103").Value = 3320.3809
$ws.Range("K103").Value = 1375.00002
$ws.Range("L103").Value = 9961.1427
$ws.Range("M103").Value = -496.0000199999999
$ws.Range("N103").Value = -11719.1427

$ws.Range("H131").Value = 22360.453
$ws.Range("I131").Value = 636.25
$ws.Range("J131").Value = 27472.03
$ws.Range("K131").Value = 1908.75
$ws.Range("L131").Value = 82416.09
$ws.Range("M131").Value = 3131.25
$ws.Range("N131").Value = -92496.09

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 6847.375
$ws.Range("I31").Value = 2111.2856
$ws.Range("K31").Value = 2111.2856
$ws.Range("M31").Value = -1819.2856

$ws.Range("H37").Value = 6847.375
$ws.Range("I37").Value = 2111.2856
$ws.Range("K37").Value = 2111.2856
$ws.Range("M37").Value = -1834.2856

$ws.Range("H102").Value = 2920.1538
$ws.Range("I102").Value = 2796.5454
$ws.Range("J102").Value = 3600
$ws.Range("K102").Value = 2796.5454
$ws.Range("L102").Value = 3600
$ws.Range("M102").Value = -1174.5454
$ws.Range("N102").Value = -6844

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 1450
$ws.Range("J46").Value = 930.7692
$ws.Range("K46").Value = 1450
$ws.Range("L46").Value = 930.7692
$ws.Range("M46").Value = -1262
$ws.Range("N46").Value = -1306.7692

$ws.Range("H61").Value = 2826.25
$ws.Range("I61").Value = 2100
$ws.Range("J61").Value = 5005
$ws.Range("K61").Value = 2100
$ws.Range("L61").Value = 5005
$ws.Range("M61").Value = -1898
$ws.Range("N61").Value = -5409

$ws.Range("H113").Value = 2826.25
$ws.Range("I113").Value = 2100
$ws.Range("J113").Value = 5005
$ws.Range("K113").Value = 2100
$ws.Range("L113").Value = 5005
$ws.Range("M113").Value = 70
$ws.Range("N113").Value = -9345

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1392.5333
$ws.Range("I113").Value = 848
$ws.Range("J113").Value = 2014.8572
$ws.Range("K113").Value = 2544
$ws.Range("L113").Value = 6044.571599999999
$ws.Range("M113").Value = -374
$ws.Range("N113").Value = -10384.5716

$ws.Range("H123").Value = 70000
$ws.Range("J123").Value = 70000
$ws.Range("L123").Value = 70000
$ws.Range("N123").Value = -79800

$ws.Range("H126").Value = 1681.3636
$ws.Range("I126").Value = 1586.8125
$ws.Range("J126").Value = 1933.5
$ws.Range("K126").Value = 4760.4375
$ws.Range("L126").Value = 5800.5
$ws.Range("M126").Value = -2290.4375
$ws.Range("N126").Value = -10740.5

$ws.Range("H132").Value = 55514.89
$ws.Range("I132").Value = 46302.41
$ws.Range("J132").Value = 69026.53
$ws.Range("K132").Value = 138907.23
$ws.Range("L132").Value = 207079.59
$ws.Range("M132").Value = -136377.23
$ws.Range("N132").Value = -212139.59

$ws.Range("H135").Value = 43248.5
$ws.Range("J135").Value = 43248.5
$ws.Range("L135").Value = 43248.5
$ws.Range("N135").Value = -53388.5
